$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 63.91118233333333
$ws.Range("H2").Value = 191.733547
$ws.Range("I2").Value = 0.4067926910433548
$ws.Range("J2").Value = 0.4067926910433549
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.040639
$ws.Range("N2").Value = 3.121917
$ws.Range("O2").Value = 0.4775850321299971
$ws.Range("P2").Value = 0.4775850321299971
$ws.Range("Q2").Value = 66.50846887217766
$ws.Range("R2").Value = 598.5762198495989
$ws.Range("S2").Value = 0.1942781004221886
$ws.Range("T2").Value = 0.1942781004221886

# Row 3
$ws.Range("G3").Value = 63.91118233333333
$ws.Range("H3").Value = 191.733547
$ws.Range("I3").Value = 0.4067926910433548
$ws.Range("J3").Value = 0.4067926910433549
$ws.Range("O3").Value = 0.2280221671432956
$ws.Range("P3").Value = 0.2280221671432955
$ws.Range("Q3").Value = 31.75435615722644
$ws.Range("R3").Value = 285.7892054150379
$ws.Range("S3").Value = 0.09275775098975884
$ws.Range("T3").Value = 0.09275775098975884

# Row 4
$ws.Range("G4").Value = 63.91118233333333
$ws.Range("H4").Value = 191.733547
$ws.Range("I4").Value = 0.4067926910433548
$ws.Range("J4").Value = 0.4067926910433549
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6414703333333333
$ws.Range("N4").Value = 1.924411
$ws.Range("O4").Value = 0.2943928007267073
$ws.Range("P4").Value = 0.2943928007267073
$ws.Range("Q4").Value = 40.99712743509077
$ws.Range("R4").Value = 368.974146915817
$ws.Range("S4").Value = 0.1197568396314074
$ws.Range("T4").Value = 0.1197568396314074

# Row 5
$ws.Range("I5").Value = 0.3656254573230189
$ws.Range("J5").Value = 0.365625457323019
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.040639
$ws.Range("N5").Value = 3.121917
$ws.Range("O5").Value = 0.4775850321299971
$ws.Range("P5").Value = 0.4775850321299971
$ws.Range("Q5").Value = 59.77784233259999
$ws.Range("R5").Value = 538.0005809933999
$ws.Range("S5").Value = 0.1746172457831589
$ws.Range("T5").Value = 0.1746172457831589

# Row 6
$ws.Range("I6").Value = 0.3656254573230189
$ws.Range("J6").Value = 0.365625457323019
$ws.Range("O6").Value = 0.2280221671432956
$ws.Range("P6").Value = 0.2280221671432955
$ws.Range("S6").Value = 0.0833707091415533
$ws.Range("T6").Value = 0.08337070914155331

# Row 7
$ws.Range("I7").Value = 0.3656254573230189
$ws.Range("J7").Value = 0.365625457323019
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6414703333333333
$ws.Range("N7").Value = 1.924411
$ws.Range("O7").Value = 0.2943928007267073
$ws.Range("P7").Value = 0.2943928007267073
$ws.Range("Q7").Value = 36.8482369458
$ws.Range("R7").Value = 331.6341325122
$ws.Range("S7").Value = 0.1076375023983067
$ws.Range("T7").Value = 0.1076375023983068

# Row 8
$ws.Range("G8").Value = 35.755375
$ws.Range("H8").Value = 107.266125
$ws.Range("I8").Value = 0.2275818516336261
$ws.Range("J8").Value = 0.2275818516336262
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.040639
$ws.Range("N8").Value = 3.121917
$ws.Range("O8").Value = 0.4775850321299971
$ws.Range("P8").Value = 0.4775850321299971
$ws.Range("Q8").Value = 37.208437684625
$ws.Range("R8").Value = 334.875939161625
$ws.Range("S8").Value = 0.1086896859246496
$ws.Range("T8").Value = 0.1086896859246496

# Row 9
$ws.Range("G9").Value = 35.755375
$ws.Range("H9").Value = 107.266125
$ws.Range("I9").Value = 0.2275818516336261
$ws.Range("J9").Value = 0.2275818516336262
$ws.Range("O9").Value = 0.2280221671432956
$ws.Range("P9").Value = 0.2280221671432955
$ws.Range("Q9").Value = 17.76510574258333
$ws.Range("R9").Value = 159.88595168325
$ws.Range("S9").Value = 0.05189370701198339
$ws.Range("T9").Value = 0.05189370701198339

# Row 10
$ws.Range("G10").Value = 35.755375
$ws.Range("H10").Value = 107.266125
$ws.Range("I10").Value = 0.2275818516336261
$ws.Range("J10").Value = 0.2275818516336262
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6414703333333333
$ws.Range("N10").Value = 1.924411
$ws.Range("O10").Value = 0.2943928007267073
$ws.Range("P10").Value = 0.2943928007267073
$ws.Range("Q10").Value = 22.93601231970833
$ws.Range("R10").Value = 206.424110877375
$ws.Range("S10").Value = 0.06699845869699317
$ws.Range("T10").Value = 0.06699845869699317
